# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-11-26 01:10:04
#
# The "Recorded By" column (G) stores a comma-separated list of identities
# that recorded/updated the session. Upstream normalized the ordering of
# a handful of known combinations so that the real actor email comes first
# and the generic "System"/"system" marker(s) trail at the end.
#
# This applies the exact same literal substitutions used upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Literal old -> new substitutions for the "Recorded By" (column G) text.
$map = @{
    "system, System, backup@backdoor.com" = "backup@backdoor.com, system, System"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
    "System, admin@admin.com"             = "admin@admin.com, System"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    $val = $cell.Value()
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
